# Generate Report for Handback
# Regenerates the localization-status report: the zh-cn / de-de handback
# rows move from "Ready for handoff" to "Handed back: in sync with en-US",
# the handback timestamps are refreshed, and the (now resolved) error
# details are cleared. Column widths are auto-resized to fit the new,
# longer Status text and the shorter (now blank) Error Detail column.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns ---
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("K2").Value = "2016-08-13 03:00:19"
$ws2.Range("P2").Value = ""

# --- de-de sheet ---
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("K2").Value = "2016-08-13 03:00:32"
$ws3.Range("P2").Value = ""

# --- Column width auto-resize (Status column widened, Error Detail narrowed) ---
# ColumnWidth is expressed in characters and snaps to the workbook's
# digit-width pixel grid, so we pick the character width that rounds to
# the closest achievable stored width to the target.
$ws1.Columns.Item(5).ColumnWidth  = 29.166666666666668   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth  = 29.166666666666668   # Overview!F (de-de status)

$ws2.Columns.Item(3).ColumnWidth  = 29.166666666666668   # zh-cn!C (Status)
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334   # zh-cn!P (Error Detail)

$ws3.Columns.Item(3).ColumnWidth  = 29.166666666666668   # de-de!C (Status)
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334   # de-de!P (Error Detail)
